$wb = $excel.ActiveWorkbook

# Both the "展览" sheet and the "全部类型" sheet contain the same first four
# data rows, and both need their F/G values (想去人数 / 最低票价) updated.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 849
    $ws.Range("G2").Value = 50

    $ws.Range("F3").Value = 4333

    $ws.Range("F4").Value = 122

    $ws.Range("F5").Value = 776
}
